$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") is a text column (values like "68.033.61" use "." as a
# thousands separator, and trailing zeros like "1.00" matter). Whenever the new
# price string would otherwise be auto-recognised as a number by Excel, prefix
# it with an apostrophe so it is stored verbatim as text, then reapply the
# "Normal" style so the forced text-prefix does not leave a stray number format.

$ws.Range("D2").Value = "68.094.21"
$ws.Range("E2").Value = "  +1.76%  "

$ws.Range("D3").Value = "3.341.56"
$ws.Range("E3").Value = "  +1.96%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'582.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.15%  "

$ws.Range("D6").Value = "'177.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.54%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").Value = "3.336.65"
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("D10").Value = "'0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.02%  "

$ws.Range("D11").Value = "'0.582"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.87%  "

$ws.Range("D12").Value = "'46.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.95%  "

$ws.Range("D13").Value = "'0.0000274"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.34%  "

$ws.Range("D14").Value = "'690.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").Value = "3.881.90"
$ws.Range("E15").Value = "  +1.90%  "

$ws.Range("D16").Value = "'8.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.14%  "

$ws.Range("D17").Value = "68.109.79"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "3.339.15"
$ws.Range("E19").Value = "  +1.58%  "

$ws.Range("D20").Value = "'17.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").Value = "'11.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.70%  "

$ws.Range("D22").Value = "'0.899"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.52%  "

$ws.Range("D23").Value = "'5.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.34%  "

$ws.Range("D24").Value = "'17.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("D25").Value = "'99.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("D26").Value = "'3.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").Value = "'9.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.36%  "

$ws.Range("D29").Value = "'33.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.81%  "

$ws.Range("D30").Value = "'8.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.79%  "

$ws.Range("D31").Value = "'7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.00%  "

$ws.Range("D32").Value = "'572.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("E33").Value = "  +2.08%  "

$ws.Range("D34").Value = "'0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.83%  "

$ws.Range("D35").Value = "3.723.23"
$ws.Range("E35").Value = "  -3.94%  "

$ws.Range("D36").Value = "'57.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.15%  "

$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").Value = "'3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").Value = "'34.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.03%  "

$ws.Range("E40").Value = "  +2.30%  "

# Rows 41 and 42 swap coin content (Fetch.AI <-> Stacks)
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.81%  "

$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.44%  "

$ws.Range("D43").Value = "0.0₃0679"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("E44").Value = "  +3.31%  "

$ws.Range("D45").Value = "'3.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("D46").Value = "'0.0408"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("D47").Value = "'2.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.91%  "

$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("E50").Value = "  -2.83%  "

$ws.Range("D51").Value = "'129.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
